$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $text)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

$sub3 = [char]0x2083

# Row 2
Set-TextValue "D2" '36.910.22'
Set-TextValue "E2" '  -0.64%  '

# Row 3
Set-TextValue "D3" '2.087.95'
Set-TextValue "E3" '  +1.76%  '

# Row 4
Set-TextValue "E4" '  +0.10%  '

# Row 5
Set-TextValue "E5" '  -0.95%  '

# Row 6
Set-TextValue "D6" '0.654'
Set-TextValue "E6" '  -1.53%  '

# Row 7
Set-TextValue "E7" '  +0.04%  '

# Row 8
Set-TextValue "D8" '55.82'
Set-TextValue "E8" '  -4.02%  '

# Row 9
Set-TextValue "D9" '59.98'
Set-TextValue "E9" '  -0.27%  '

# Row 10
Set-TextValue "D10" '0.369'
Set-TextValue "E10" '  -3.46%  '

# Row 11
Set-TextValue "D11" '0.0769'
Set-TextValue "E11" '  -1.24%  '

# Row 12
Set-TextValue "E12" '  +1.39%  '

# Row 13
Set-TextValue "D13" '15.14'
Set-TextValue "E13" '  -4.65%  '

# Row 14
Set-TextValue "D14" '0.890'
Set-TextValue "E14" '  +6.61%  '

# Row 15
Set-TextValue "D15" '2.401.06'
Set-TextValue "E15" '  +2.02%  '

# Row 16
Set-TextValue "B16" 'WrappedEther'
Set-TextValue "C16" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D16" '2.222.28'
Set-TextValue "E16" '  +7.85%  '

# Row 17
Set-TextValue "B17" 'Polkadot'
Set-TextValue "C17" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D17" '5.54'
Set-TextValue "E17" '  -2.77%  '

# Row 18
Set-TextValue "D18" '36.875.70'
Set-TextValue "E18" '  -0.90%  '

# Row 19
Set-TextValue "D19" '17.43'
Set-TextValue "E19" '  -3.73%  '

# Row 20
Set-TextValue "D20" '73.22'
Set-TextValue "E20" '  -2.07%  '

# Row 21
Set-TextValue "D21" ("0.0{0}0886" -f $sub3)
Set-TextValue "E21" '  -1.14%  '

# Row 22
Set-TextValue "D22" '5.50'
Set-TextValue "E22" '  +2.89%  '

# Row 23
Set-TextValue "D23" '237.63'
Set-TextValue "E23" '  +0.27%  '

# Row 24
Set-TextValue "D24" '0.999'
Set-TextValue "E24" '  -0.06%  '

# Row 25
Set-TextValue "D25" '2.42'
Set-TextValue "E25" '  -1.41%  '

# Row 26
Set-TextValue "D26" '9.89'
Set-TextValue "E26" '  +4.91%  '

# Row 27
Set-TextValue "E27" '  +0.09%  '

# Row 28
Set-TextValue "D28" '168.62'
Set-TextValue "E28" '  -0.39%  '

# Row 29
Set-TextValue "D29" '20.72'
Set-TextValue "E29" '  +3.43%  '

# Row 30
Set-TextValue "D30" '5.45'
Set-TextValue "E30" '  +13.70%  '

# Row 31
Set-TextValue "D31" '0.123'
Set-TextValue "E31" '  -0.37%  '

# Row 32
Set-TextValue "D32" '1.19'
Set-TextValue "E32" '  +6.35%  '

# Row 33
Set-TextValue "E33" '  +5.13%  '

# Row 34
Set-TextValue "D34" '0.0612'
Set-TextValue "E34" '  -1.00%  '

# Row 35
Set-TextValue "D35" '2.39'
Set-TextValue "E35" '  +5.97%  '

# Row 36
Set-TextValue "E36" '  +0.19%  '

# Row 37
Set-TextValue "E37" '  +4.28%  '

# Row 38
Set-TextValue "D38" '0.0842'
Set-TextValue "E38" '  -6.11%  '

# Row 39
Set-TextValue "E39" '  -3.42%  '

# Row 40
Set-TextValue "D40" '1.16'
Set-TextValue "E40" '  +2.09%  '

# Row 41
Set-TextValue "E41" '  -0.11%  '

# Row 42
Set-TextValue "D42" '4.89'
Set-TextValue "E42" '  -6.17%  '

# Row 43
Set-TextValue "D43" '0.0952'
Set-TextValue "E43" '  -6.70%  '

# Row 44
Set-TextValue "D44" '97.01'
Set-TextValue "E44" '  +1.32%  '

# Row 45
Set-TextValue "D45" '2.86'
Set-TextValue "E45" '  -11.76%  '

# Row 46
Set-TextValue "D46" '16.12'
Set-TextValue "E46" '  -6.09%  '

# Row 47
Set-TextValue "D47" '1.354.15'
Set-TextValue "E47" '  +6.16%  '

# Row 48
Set-TextValue "D48" '2.46'
Set-TextValue "E48" '  +0.50%  '

# Row 49
Set-TextValue "B49" 'MXToken'
Set-TextValue "C49" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D49" '2.91'
Set-TextValue "E49" '  +1.49%  '

# Row 50
Set-TextValue "B50" 'FraxShare'
Set-TextValue "C50" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D50" '7.01'
Set-TextValue "E50" '  +2.89%  '

# Row 51
Set-TextValue "D51" '2.284.00'
Set-TextValue "E51" '  +2.04%  '
